$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename model headers (row 1)
$ws.Range("C1").Value = "Modelo 1"
$ws.Range("D1").Value = "Modelo 2"
$ws.Range("E1").Value = "Modelo 3"

# Update (Intercept) estimates for Model 2 and Model 3 (row 2)
$ws.Range("D2").Value = "-177063223.674***"
$ws.Range("E2").Value = "-69043404.727***"

# Update (Intercept) std errors for Model 2 and Model 3 (row 3)
$ws.Range("D3").Value = "[7052644.424]"
$ws.Range("E3").Value = "[8127293.665]"

# Rename "tipo" variable to "as.factor(property_type)Casa" (row 10)
$ws.Range("A10").Value = "as.factor(property_type)Casa"

# Flip sign of the estimates for that variable in Model 2 and Model 3
$ws.Range("D10").Value = "-217848331.263***"
$ws.Range("E10").Value = "-173754264.622***"
